$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (R1) - uptime changed
$ws.Range("G3").Value = "5:36:00"

# Row 4: was R3 -> now R2 (uptime also refreshed)
$ws.Range("A4").Value = "R2"
$ws.Range("B4").Value = "R2.automation.local"
$ws.Range("G4").Value = "5:42:22"

# Row 5: was SW1 -> now R3, re-using the R1 device profile (Cisco 7206VXR router), new uptime
$ws.Range("A5").Value = "R3"
$ws.Range("B5").Value = "R3.automation.local"
$ws.Range("C5").Value = "Cisco"
$ws.Range("D5").Value = "7206VXR"

# E5 ("4279256517") is numeric-looking text, not a number - force text storage
# (NumberFormat "@" stops auto-coercion), then copy E3's format back over so the
# cell keeps the sheet's normal (centered) style instead of a new "text" style.
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4279256517"
$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F5").Value = "Version 15.2(4)S5"
$ws.Range("G5").Value = "5:34:00"
$ws.Range("H5").Value = "OK"

# Row 6 (new): MLS1 - takes over the old SW1 device profile (Cisco IOSv), new serial + uptime
# Match the sheet's existing centered cell style for the new row (copy format
# from row 2, which already uses that style, instead of re-deriving alignment
# and minting a brand new style entry).
$ws.Range("A2:H2").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A6").Value = "MLS1"
$ws.Range("B6").Value = "MLS1.automation.local"
$ws.Range("C6").Value = "Cisco"
$ws.Range("D6").Value = "IOSv"
$ws.Range("E6").Value = "9LD1YQMD0KM"
$ws.Range("F6").Value = "Version 15.2(4.0.55)E"
$ws.Range("G6").Value = "5:42:00"
$ws.Range("H6").Value = "OK"
